$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O1:Q1").EntireColumn.Delete()
$ws.Range("A3:A10").WrapText = $true
$ws.Range("A1:N10").Select()
$excel.ActiveWindow.Zoom = 85
Write-Output "done"
